$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: only B2 changes slightly
$ws.Range("B2").Value = 0.1014141115953233

# Row 3 - RandomForestRegressor: B3, C3, D3 change
$ws.Range("B3").Value = 0.01783216377055935
$ws.Range("C3").Value = 0.0166114838558714
$ws.Range("D3").Value = 0.01760415828564628

# Row 4 - model name changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01648238219824683
$ws.Range("C4").Value = 0.01488446158577947
$ws.Range("D4").Value = 0.01479089809016145

# Row 5 - model name changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01425077639215865
$ws.Range("C5").Value = 0.0146282420983529
$ws.Range("D5").Value = 0.01495697851874981
